$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldTop = '"Part Number - Can be found on the top right position of the page"'
$newCenter = '"Part Number - Can be found on the center right position of the page"'

$oldDuxford = '"Duxford Range Part Number Description Dimensions Power Lumens Colour Temp. - Can be found on the bottom right position of the page"'
$newMulti = '"Multi-Wattage Tri-Colour and Single Colour 4000K Retrofit Gear Trays - Can be found on the middle right position of the page"'

for ($r = 2; $r -le 13; $r++) {
    $ws.Cells.Item($r, 7).Value = $newCenter
}

for ($r = 14; $r -le 25; $r++) {
    $ws.Cells.Item($r, 7).Value = $newMulti
}
